$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: new header date cell BB1, copy formatting (style) from BA1 then set its value
$ws.Range("BA1").Copy()
$ws.Range("BB1").PasteSpecial(-4122)
$ws.Range("BB1").Value = 45986

# Rows 3-18: new BB column value equal to the existing BA column value for that row
$ws.Range("BB3").Value = -3.712113773657133
$ws.Range("BB4").Value = 2.363997455008704
$ws.Range("BB5").Value = 6.087058041562399
$ws.Range("BB6").Value = 3.831642495756782
$ws.Range("BB7").Value = -0.492746714819392
$ws.Range("BB8").Value = 4.135644887697842
$ws.Range("BB9").Value = -0.8017097776176652
$ws.Range("BB10").Value = 2.496356089727558
$ws.Range("BB11").Value = 2.955681454251202
$ws.Range("BB12").Value = 2.115640452469392
$ws.Range("BB13").Value = 4.186780469526941
$ws.Range("BB14").Value = 1.25858038212967
$ws.Range("BB15").Value = 1.362828756617751
$ws.Range("BB16").Value = -2.151676542405401
$ws.Range("BB17").Value = -3.961484026309636
$ws.Range("BB18").Value = -2.899559879361435

# Rows 19-21: new BB column values (forecast continuation, differ from BA column)
$ws.Range("BB19").Value = -0.7200474048664085
$ws.Range("BB20").Value = -0.2284091334091687
$ws.Range("BB21").Value = 0.534104719104489

# Row 22 intentionally left without a BB cell, matching the source diff
